$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Paragraph 1 ("Resolved"): the text is re-typed as two pieces
# ("Res" then "olved") split around the existing _GoBack bookmark.
# ---------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$titleRange = $d.Range($p1.Range.Start, $p1.Range.Start + 8)
$titleRange.Text = "Res"
$bm = $d.Bookmarks.Item("_GoBack")
$afterBookmark = $d.Bookmarks.Item("_GoBack").End
$tail = $d.Range($afterBookmark, $afterBookmark)
$tail.InsertAfter("olved")

# ---------------------------------------------------------------
# Paragraph 2 (was empty / centered / sz 68): becomes note #1,
# left-aligned, Times New Roman 12pt bold.
# ---------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$p2.Alignment = 0
$seed2 = $p2.Range
$seed2.Font.NameBi = "Times New Roman"
$seed2.Font.SizeBi = 12
$seed2.InsertAfter("1.Added ajax to all the pages that needs it, except universitydownload2.php. This fix the issue with the back button and the navigation is much smoother.")
$full2 = $p2.Range
$full2.Font.Name = "Times New Roman"
$full2.Font.NameBi = "Times New Roman"
$full2.Font.Size = 12
$full2.Font.SizeBi = 12
$full2.Font.Bold = $true
$full2.Font.Underline = 0

# ---------------------------------------------------------------
# Paragraph 3 (was empty / sz 48): becomes note #2, Times New
# Roman 12pt bold.
# ---------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$seed3 = $p3.Range
$seed3.Font.NameBi = "Times New Roman"
$seed3.Font.SizeBi = 12
$seed3.InsertAfter("2.Changed the button description is facultySearch.php and universitySearch.php from saying their Id number to saying " + [char]0x2018 + "Go" + [char]0x2019 + ".")
$full3 = $p3.Range
$full3.Font.Name = "Times New Roman"
$full3.Font.NameBi = "Times New Roman"
$full3.Font.Size = 12
$full3.Font.SizeBi = 12
$full3.Font.Bold = $true

# ---------------------------------------------------------------
# New paragraph 4: note #3. Inherit paragraph-mark formatting by
# splitting off of paragraph 3, then fill in its text + formatting.
# ---------------------------------------------------------------
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Item(4)
$seed4 = $p4.Range
$seed4.Font.NameBi = "Times New Roman"
$seed4.Font.SizeBi = 12
$seed4.InsertAfter("3. Added a ref to all the change setting pages such as changeEmail to go back to the profile page. Except changepassword, because their a way to get to that page if you forgot your password.")
$full4 = $p4.Range
$full4.Font.Name = "Times New Roman"
$full4.Font.NameBi = "Times New Roman"
$full4.Font.Size = 12
$full4.Font.SizeBi = 12
$full4.Font.Bold = $true

# ---------------------------------------------------------------
# New paragraph 5: note #4.
# ---------------------------------------------------------------
$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs.Item(5)
$seed5 = $p5.Range
$seed5.Font.NameBi = "Times New Roman"
$seed5.Font.SizeBi = 12
$seed5.InsertAfter("4. Add a dropdown bar with all the university in the database for the FirstPage.php. Also added a separate page to add a university if it not in the database.")
$full5 = $p5.Range
$full5.Font.Name = "Times New Roman"
$full5.Font.NameBi = "Times New Roman"
$full5.Font.Size = 12
$full5.Font.SizeBi = 12
$full5.Font.Bold = $true

Write-Host "Paragraphs:" $d.Paragraphs.Count
